$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows that were dropped from the dataset (row 28 first so row 26 index stays valid)
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Apply updated / imputed / cleared cell values
$ws.Range("D2").Value = -13.5
$ws.Range("F2").Value = 18.03
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("E5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("F13").Value = 17.1
$ws.Range("D14").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("D22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F25").Value = 16.6
$ws.Range("E27").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("E29").Value = -6.8
$ws.Range("B30").Value = -19.7
$ws.Range("D31").Value = -13.7
$ws.Range("F31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("F32").Value = 17.39
$ws.Range("D33").Value = -14.1
